# Spells.xlsx - "Made starter decks for three of the mythologys."
#
# Adds starter-deck spell lists (columns H:Q) for the WATER, WAR, LOVE,
# KING, DEATH and MESSENGER mythologies to rows 10-18 of Sheet1, and
# updates the sheet selection/view accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("H10").Value = "new"

# --- New starter-deck cells, entered in the same order the author typed
#     them (this keeps the shared-string table identical to the source). ---
$ws.Range("I12").Value = "Adrenaline Rush"
$ws.Range("I13").Value = "Cure"
$ws.Range("J13").Value = "Circle of Protection"
$ws.Range("K13").Value = "Guardian Circle"
$ws.Range("I17").Value = "God's name + Aura"
$ws.Range("I16").Value = "Familiar"
$ws.Range("I15").Value = "Ash Shell"
$ws.Range("I14").Value = "Cloud's Call"
$ws.Range("I11").Value = "Rising Tide"
$ws.Range("J15").Value = "Styx burn"
$ws.Range("L13").Value = "Mesmerize"
$ws.Range("J16").Value = "You've got mail"
$ws.Range("J12").Value = "War Chant/Song"
$ws.Range("M13").Value = "Stone Heart"
$ws.Range("I18").Value = "light"
$ws.Range("J11").Value = "Freeze"
$ws.Range("N13").Value = "Pacify"
$ws.Range("O13").Value = "Seduction"
$ws.Range("J14").Value = "Aristocracy"
$ws.Range("K14").Value = "Overlord"
$ws.Range("L14").Value = "Conquering Hero"
$ws.Range("M14").Value = "Vanquisher"
$ws.Range("N14").Value = "Authority"
$ws.Range("P13").Value = "Awe"
$ws.Range("O14").Value = "Coronation"
$ws.Range("P14").Value = "pull rank"
$ws.Range("Q14").Value = "Allegiance"
$ws.Range("K11").Value = "calm"
$ws.Range("L11").Value = "rain"
$ws.Range("M11").Value = "typhoon"
$ws.Range("N11").Value = "WaterWell"
$ws.Range("O11").Value = "Soak"
$ws.Range("P11").Value = "steam"
$ws.Range("K15").Value = "after life"
$ws.Range("L15").Value = "Blood Shed"
$ws.Range("M15").Value = "Eternal rest"

# --- Mythology-name labels in column H (reuse existing shared strings). ---
$ws.Range("H11").Value = "WATER"
$ws.Range("H12").Value = "WAR"
$ws.Range("H13").Value = "LOVE"
$ws.Range("H14").Value = "KING"
$ws.Range("H15").Value = "DEATH"
$ws.Range("H16").Value = "MESSENGER"

# --- View / selection changes ---
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("N15").Select()
